# Update "Seasonality Index" (column L) values on the "Forecast Comparison" sheet
# to reflect the refreshed inventory-file read.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$values = @{
    2  = 0.85
    3  = 0.99
    4  = 1.02
    5  = 1.08
    6  = 0.9399999999999999
    7  = 0.9
    8  = 0.8100000000000001
    9  = 1.07
    10 = 0.98
    11 = 0.9
    12 = 1
    13 = 0.91
    14 = 1.09
    15 = 1.13
    16 = 0.92
    17 = 1.11
}

foreach ($row in $values.Keys) {
    $ws.Range("L$row").Value = $values[$row]
}
